$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The commit rotates the species-observation data stored in rows 75, 76
# and 78 of the "Artfynd" sheet (row 77 keeps its own data and is not
# touched). The location/date/survey metadata columns (C,K,P,T,U,V,W,Y,
# AA,AD,AE,AG,AT,AY) are identical across rows 75/76/78 already, so only
# the columns that actually differ between the three rows need to be
# rewritten:
#
#   new row75  <-  old row78   (Knärot / Goodyera repens, Kamilla Andersson)
#   new row76  <-  old row75   (Fläcknycklar / Dactylorhiza maculata)
#   new row78  <-  old row76   (Talltita / Poecile montanus)
#
# Column I holds numeric-looking values that must stay text (e.g. "15"),
# and row75 gains a J cell ("plantor/tuvor") that row78 loses.

function Set-TextValue($range, [string]$text) {
    # Force a numeric-looking string to be stored as text, the same way
    # the original file stores it (quantity "15", "7", ... as text), then
    # drop the temporary Text number-format so no stray style lingers.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- New row 75 values (= old row 78) ---
$ws.Range("A75").Value = 111949575
$ws.Range("B75").Value = 96348
$ws.Range("D75").Value = "VU"
$ws.Range("E75").Value = 220787
$ws.Range("F75").Value = "Knärot"
$ws.Range("G75").Value = "Goodyera repens"
$ws.Range("H75").Value = "(L.) R. Br."
Set-TextValue $ws.Range("I75") "15"
Set-TextValue $ws.Range("J75") "plantor/tuvor"
$ws.Range("Q75").Value = 580471.3517951096
$ws.Range("R75").Value = 7053333.257918903
$ws.Range("S75").Value = 1
$ws.Range("Z75").Value = "19:05"
$ws.Range("AB75").Value = "19:05"
$ws.Range("AW75").Value = "Kamilla Andersson"
$ws.Range("AX75").Value = "Kamilla Andersson"

# --- New row 76 values (= old row 75) ---
$ws.Range("A76").Value = 111949317
$ws.Range("B76").Value = 96265
$ws.Range("D76").Value = "LC"
$ws.Range("E76").Value = 219790
$ws.Range("F76").Value = "Fläcknycklar"
$ws.Range("G76").Value = "Dactylorhiza maculata"
$ws.Range("H76").Value = "(L.) Soó"
$ws.Range("I76").ClearContents()
$ws.Range("Q76").Value = 580500.003505226
$ws.Range("R76").Value = 7053328.641698814
$ws.Range("S76").Value = 2
$ws.Range("Z76").Value = "18:54"
$ws.Range("AB76").Value = "18:54"
$ws.Range("AW76").Value = "Kim Hultgren"
$ws.Range("AX76").Value = "Kim Hultgren"

# --- New row 78 values (= old row 76) ---
$ws.Range("A78").Value = 111950184
$ws.Range("B78").Value = 56543
$ws.Range("D78").Value = "NT"
$ws.Range("E78").Value = 103021
$ws.Range("F78").Value = "Talltita"
$ws.Range("G78").Value = "Poecile montanus"
$ws.Range("H78").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I78").ClearContents()
$ws.Range("J78").ClearContents()
$ws.Range("Q78").Value = 580446.7330953531
$ws.Range("R78").Value = 7053301.910512885
$ws.Range("S78").Value = 10
$ws.Range("Z78").Value = "19:37"
$ws.Range("AB78").Value = "19:37"
$ws.Range("AW78").Value = "Kim Hultgren"
$ws.Range("AX78").Value = "Kim Hultgren"
